$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Physiology")
$ws.Range("A1").Value = "test"
